# ABD dip hourly update
# The ticker list in column A shifted: 12 symbols were removed and 4 new
# symbols were inserted at specific points, with everything else keeping
# its relative order. Apply the row deletions/insertions from the bottom
# of the sheet upward so that earlier (still-to-be-processed) row numbers
# are not disturbed by the edits already made below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Deletions (rows removed from the list) ---------------------------
# row 434 -> DDS
$ws.Rows(434).Delete()
# row 423 -> JAZZ
$ws.Rows(423).Delete()
# row 396 -> CPT
$ws.Rows(396).Delete()

# --- Insertion: NVMI (was not present before) --------------------------
$ws.Rows(355).Insert()
$ws.Range("A355").Value = "NVMI"

# --- Deletions continued ------------------------------------------------
# row 325 -> ONON
$ws.Rows(325).Delete()
# row 322 -> FFIV
$ws.Rows(322).Delete()
# row 305 -> BSAC
$ws.Rows(305).Delete()
# row 287 -> FTV
$ws.Rows(287).Delete()
# row 235 -> FLEX
$ws.Rows(235).Delete()

# --- Insertion: DLR-PK ---------------------------------------------------
$ws.Rows(230).Insert()
$ws.Range("A230").Value = "DLR-PK"

# --- Deletions continued ------------------------------------------------
# row 227 -> WRB
$ws.Rows(227).Delete()
# row 209 -> HPE
$ws.Rows(209).Delete()

# --- Insertion: MS-PA -----------------------------------------------------
$ws.Rows(112).Insert()
$ws.Range("A112").Value = "MS-PA"

# --- Insertion: RCL ---------------------------------------------------------
$ws.Rows(100).Insert()
$ws.Range("A100").Value = "RCL"

# --- Deletions continued ------------------------------------------------
# row 88 -> PNC
$ws.Rows(88).Delete()
# row 83 -> NKE
$ws.Rows(83).Delete()
